$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New daily series rows: 09-10-2021 .. 14-10-2021 (rows 283-288)
$dates = @("09-10-2021", "10-10-2021", "11-10-2021", "12-10-2021", "13-10-2021", "14-10-2021")

$data = @(
    @(2083, 2598, 15028, 2623, 4542, 8726),
    @(2083, 2598, 15028, 2623, 4542, 8726),
    @(2083, 2598, 15028, 2623, 4542, 8726),
    @(2066, 2578, 14908, 2602, 4506, 8657),
    @(2048, 2555, 14775, 2579, 4466, 8579),
    @(2066, 2578, 14908, 2602, 4506, 8657)
)

$startRow = 283

# Force column A to be read as text so the dd-mm-yyyy strings are stored
# as literal text (matching the existing "Serie" column) instead of being
# auto-converted into date serial numbers.
$dateRange = $ws.Range("A" + $startRow + ":A" + ($startRow + $dates.Length - 1))
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i

    $ws.Cells.Item($r, 1).Value = $dates[$i]

    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, 2 + $j).Value = $row[$j]
    }
}

# Restore the default (unstyled) look for column A on the new rows, same as
# the rest of the "Serie" column.
$dateRange.Style = "Normal"
